$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 13499.333
$ws.Cells.Item(21, 9).Value = 500
$ws.Cells.Item(21, 11).Value = 500
$ws.Cells.Item(21, 13).Value = -32
$ws.Cells.Item(23, 8).Value = 13499.333
$ws.Cells.Item(23, 9).Value = 500
$ws.Cells.Item(23, 11).Value = 500
$ws.Cells.Item(23, 13).Value = -266
$ws.Cells.Item(32, 8).Value = 1209.125
$ws.Cells.Item(32, 9).Value = 1500.5
$ws.Cells.Item(32, 10).Value = 1112
$ws.Cells.Item(32, 11).Value = 1500.5
$ws.Cells.Item(32, 12).Value = 1112
$ws.Cells.Item(32, 13).Value = -1174.5
$ws.Cells.Item(32, 14).Value = -1764
$ws.Cells.Item(52, 8).Value = 3169.9565
$ws.Cells.Item(52, 10).Value = 3223.1365
$ws.Cells.Item(52, 12).Value = 9669.4095
$ws.Cells.Item(52, 14).Value = -9989.4095
$ws.Cells.Item(70, 8).Value = 1623.7
$ws.Cells.Item(70, 9).Value = 1278.3334
$ws.Cells.Item(70, 10).Value = 1771.7142
$ws.Cells.Item(70, 11).Value = 3835.0002
$ws.Cells.Item(70, 12).Value = 5315.142599999999
$ws.Cells.Item(70, 13).Value = -3565.0002
$ws.Cells.Item(70, 14).Value = -5855.142599999999
$ws.Cells.Item(73, 8).Value = 1623.7
$ws.Cells.Item(73, 9).Value = 1278.3334
$ws.Cells.Item(73, 10).Value = 1771.7142
$ws.Cells.Item(73, 11).Value = 3835.0002
$ws.Cells.Item(73, 12).Value = 5315.142599999999
$ws.Cells.Item(73, 13).Value = -2899.0002
$ws.Cells.Item(73, 14).Value = -7187.142599999999
$ws.Cells.Item(113, 8).Value = 1930.5714
$ws.Cells.Item(113, 9).Value = 1878.75
$ws.Cells.Item(113, 10).Value = 1999.6666
$ws.Cells.Item(113, 11).Value = 1878.75
$ws.Cells.Item(113, 12).Value = 1999.6666
$ws.Cells.Item(113, 13).Value = 1375.25
$ws.Cells.Item(113, 14).Value = -8507.6666
$ws.Cells.Item(127, 8).Value = 1053.9333
$ws.Cells.Item(127, 9).Value = 504.5
$ws.Cells.Item(127, 10).Value = 1420.2222
$ws.Cells.Item(127, 11).Value = 1513.5
$ws.Cells.Item(127, 12).Value = 4260.6666
$ws.Cells.Item(127, 13).Value = 3446.5
$ws.Cells.Item(127, 14).Value = -14180.6666
$ws.Cells.Item(129, 8).Value = 1082.4314
$ws.Cells.Item(129, 10).Value = 1202.2927
$ws.Cells.Item(129, 12).Value = 3606.8781
$ws.Cells.Item(129, 14).Value = -13606.8781
$ws.Cells.Item(132, 8).Value = 4244.086
$ws.Cells.Item(132, 9).Value = 3898.303
$ws.Cells.Item(132, 11).Value = 11694.909
$ws.Cells.Item(132, 13).Value = -9164.909
$ws.Cells.Item(138, 8).Value = 2029.5671
$ws.Cells.Item(138, 9).Value = 2149
$ws.Cells.Item(138, 10).Value = 1975.0435
$ws.Cells.Item(138, 11).Value = 6447
$ws.Cells.Item(138, 12).Value = 5925.1305
$ws.Cells.Item(138, 13).Value = -1307
$ws.Cells.Item(138, 14).Value = -16205.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1310.625
$ws.Cells.Item(2, 9).Value = 1347.4286
$ws.Cells.Item(2, 10).Value = 1053
$ws.Cells.Item(2, 11).Value = 1347.4286
$ws.Cells.Item(2, 12).Value = 1053
$ws.Cells.Item(2, 13).Value = -1234.4286
$ws.Cells.Item(2, 14).Value = -1279
$ws.Cells.Item(45, 8).Value = 2926.6956
$ws.Cells.Item(45, 9).Value = 2250
$ws.Cells.Item(45, 10).Value = 3664.9092
$ws.Cells.Item(45, 11).Value = 2250
$ws.Cells.Item(45, 12).Value = 3664.9092
$ws.Cells.Item(45, 13).Value = -1873
$ws.Cells.Item(45, 14).Value = -4418.9092
$ws.Cells.Item(61, 8).Value = 4092.889
$ws.Cells.Item(61, 9).Value = 3807.2
$ws.Cells.Item(61, 10).Value = 4450
$ws.Cells.Item(61, 11).Value = 3807.2
$ws.Cells.Item(61, 12).Value = 4450
$ws.Cells.Item(61, 13).Value = -3595.2
$ws.Cells.Item(61, 14).Value = -4874
$ws.Cells.Item(92, 8).Value = 54589.285
$ws.Cells.Item(92, 10).Value = 54589.285
$ws.Cells.Item(92, 12).Value = 54589.285
$ws.Cells.Item(92, 14).Value = -59581.285
$ws.Cells.Item(116, 8).Value = 1310.625
$ws.Cells.Item(116, 9).Value = 1347.4286
$ws.Cells.Item(116, 10).Value = 1053
$ws.Cells.Item(116, 11).Value = 1347.4286
$ws.Cells.Item(116, 12).Value = 1053
$ws.Cells.Item(116, 13).Value = 946.5714
$ws.Cells.Item(116, 14).Value = -5641
$ws.Cells.Item(136, 8).Value = 4092.889
$ws.Cells.Item(136, 9).Value = 3807.2
$ws.Cells.Item(136, 10).Value = 4450
$ws.Cells.Item(136, 11).Value = 11421.6
$ws.Cells.Item(136, 12).Value = 13350
$ws.Cells.Item(136, 13).Value = -8871.599999999999
$ws.Cells.Item(136, 14).Value = -18450

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1310.625
$ws.Cells.Item(3, 9).Value = 1347.4286
$ws.Cells.Item(3, 10).Value = 1053
$ws.Cells.Item(3, 11).Value = 1347.4286
$ws.Cells.Item(3, 12).Value = 1053
$ws.Cells.Item(3, 13).Value = -1233.4286
$ws.Cells.Item(3, 14).Value = -1281
$ws.Cells.Item(20, 8).Value = 2631.3684
$ws.Cells.Item(20, 9).Value = 2624.6667
$ws.Cells.Item(20, 10).Value = 2642.8572
$ws.Cells.Item(20, 11).Value = 2624.6667
$ws.Cells.Item(20, 12).Value = 2642.8572
$ws.Cells.Item(20, 13).Value = -2377.6667
$ws.Cells.Item(20, 14).Value = -3136.8572
$ws.Cells.Item(100, 8).Value = 53598.75
$ws.Cells.Item(100, 10).Value = 53598.75
$ws.Cells.Item(100, 12).Value = 53598.75
$ws.Cells.Item(100, 14).Value = -55762.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(92, 8).Value = 41392.5
$ws.Cells.Item(92, 10).Value = 41392.5
$ws.Cells.Item(92, 12).Value = 41392.5
$ws.Cells.Item(92, 14).Value = -46384.5
$ws.Cells.Item(93, 8).Value = 23567.834
$ws.Cells.Item(93, 9).Value = 19135.666
$ws.Cells.Item(93, 10).Value = 28000
$ws.Cells.Item(93, 11).Value = 19135.666
$ws.Cells.Item(93, 12).Value = 28000
$ws.Cells.Item(93, 13).Value = -17263.666
$ws.Cells.Item(93, 14).Value = -31744
$ws.Cells.Item(96, 8).Value = 35671.43
$ws.Cells.Item(96, 10).Value = 35671.43
$ws.Cells.Item(96, 12).Value = 35671.43
$ws.Cells.Item(96, 14).Value = -41163.43
$ws.Cells.Item(100, 8).Value = 49995
$ws.Cells.Item(100, 10).Value = 49995
$ws.Cells.Item(100, 12).Value = 49995
$ws.Cells.Item(100, 14).Value = -52159
$ws.Cells.Item(112, 8).Value = 34749.75
$ws.Cells.Item(112, 10).Value = 34749.75
$ws.Cells.Item(112, 12).Value = 34749.75
$ws.Cells.Item(112, 14).Value = -37703.75
$ws.Cells.Item(119, 8).Value = 50000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 14).Value = -59676
$ws.Cells.Item(134, 8).Value = 2307.6
$ws.Cells.Item(134, 9).Value = 2881
$ws.Cells.Item(134, 10).Value = 14
$ws.Cells.Item(134, 11).Value = 8643
$ws.Cells.Item(134, 12).Value = 42
$ws.Cells.Item(134, 13).Value = -6108
$ws.Cells.Item(134, 14).Value = -5112

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1085.5
$ws.Cells.Item(5, 9).Value = 1035.8334
$ws.Cells.Item(5, 10).Value = 1234.5
$ws.Cells.Item(5, 11).Value = 3107.5002
$ws.Cells.Item(5, 12).Value = 3703.5
$ws.Cells.Item(5, 13).Value = -2995.5002
$ws.Cells.Item(5, 14).Value = -3927.5
$ws.Cells.Item(39, 8).Value = 1583.6666
$ws.Cells.Item(39, 10).Value = 1583.6666
$ws.Cells.Item(39, 12).Value = 4750.9998
$ws.Cells.Item(39, 14).Value = -5338.9998
$ws.Cells.Item(110, 8).Value = 13737.909
$ws.Cells.Item(110, 9).Value = 5009
$ws.Cells.Item(110, 11).Value = 15027
$ws.Cells.Item(110, 13).Value = -10937
$ws.Cells.Item(131, 8).Value = 1019.125
$ws.Cells.Item(131, 10).Value = 1120.8214
$ws.Cells.Item(131, 12).Value = 3362.4642
$ws.Cells.Item(131, 14).Value = -13442.4642
$ws.Cells.Item(134, 8).Value = 5391.6113
$ws.Cells.Item(134, 9).Value = 3881.125
$ws.Cells.Item(134, 10).Value = 6600
$ws.Cells.Item(134, 11).Value = 11643.375
$ws.Cells.Item(134, 12).Value = 19800
$ws.Cells.Item(134, 13).Value = -6573.375
$ws.Cells.Item(134, 14).Value = -29940
$ws.Cells.Item(135, 8).Value = 1085.5
$ws.Cells.Item(135, 9).Value = 1035.8334
$ws.Cells.Item(135, 10).Value = 1234.5
$ws.Cells.Item(135, 11).Value = 9322.500599999999
$ws.Cells.Item(135, 12).Value = 11110.5
$ws.Cells.Item(135, 13).Value = -6787.500599999999
$ws.Cells.Item(135, 14).Value = -16180.5
$ws.Cells.Item(139, 8).Value = 3365.7778
$ws.Cells.Item(139, 9).Value = 2573.625
$ws.Cells.Item(139, 10).Value = 3999.5
$ws.Cells.Item(139, 11).Value = 7720.875
$ws.Cells.Item(139, 12).Value = 11998.5
$ws.Cells.Item(139, 13).Value = -2580.875
$ws.Cells.Item(139, 14).Value = -22278.5
$ws.Cells.Item(141, 8).Value = 4800.7617
$ws.Cells.Item(141, 9).Value = 1757.6666
$ws.Cells.Item(141, 10).Value = 7083.0835
$ws.Cells.Item(141, 11).Value = 5272.9998
$ws.Cells.Item(141, 12).Value = 21249.2505
$ws.Cells.Item(141, 13).Value = -92.9997999999996
$ws.Cells.Item(141, 14).Value = -31609.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(104, 8).Value = 24547
$ws.Cells.Item(104, 10).Value = 24547
$ws.Cells.Item(104, 12).Value = 24547
$ws.Cells.Item(104, 14).Value = -31535
$ws.Cells.Item(130, 8).Value = 39795
$ws.Cells.Item(130, 10).Value = 39795
$ws.Cells.Item(130, 12).Value = 39795
$ws.Cells.Item(130, 14).Value = -49835
$ws.Cells.Item(140, 8).Value = 76476
$ws.Cells.Item(140, 10).Value = 76476
$ws.Cells.Item(140, 12).Value = 76476
$ws.Cells.Item(140, 14).Value = -86836

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 56710.668
$ws.Cells.Item(42, 9).Value = 30034
$ws.Cells.Item(42, 11).Value = 30034
$ws.Cells.Item(42, 13).Value = -29656
$ws.Cells.Item(46, 8).Value = 51546.855
$ws.Cells.Item(46, 10).Value = 51546.855
$ws.Cells.Item(46, 12).Value = 51546.855
$ws.Cells.Item(46, 14).Value = -52008.855
$ws.Cells.Item(112, 8).Value = 49990
$ws.Cells.Item(112, 10).Value = 49990
$ws.Cells.Item(112, 12).Value = 49990
$ws.Cells.Item(112, 14).Value = -52944
$ws.Cells.Item(125, 8).Value = 49199.668
$ws.Cells.Item(125, 10).Value = 49199.668
$ws.Cells.Item(125, 12).Value = 49199.668
$ws.Cells.Item(125, 14).Value = -59039.668
$ws.Cells.Item(134, 8).Value = 51546.855
$ws.Cells.Item(134, 10).Value = 51546.855
$ws.Cells.Item(134, 12).Value = 154640.565
$ws.Cells.Item(134, 14).Value = -159710.565

Write-Host "Applied all Anima_Profits updates"